$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "As we work harder" -> "As we work a little harder"
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("s we work harder", $true, $false, $false, $false, $false, $true, 1, $false, "s we work a little harder", 2) | Out-Null

# ------------------------------------------------------------------
# 2) "...to a point where we start to accumulate lactate in the legs,
#    this point is known as" -> "...to a point where lactate levels
#    increase above the baseline (typically 1mmol above), this point
#    is known as"
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("we get to a point where we start to accumulate lactate in the legs, this point is known as", $true, $false, $false, $false, $false, $true, 1, $false, "we get to a point where lactate levels increase above the baseline (typically 1mmol above), this point is known as", 2) | Out-Null

# Move the _GoBack bookmark so that it now sits right before "1mmol above)"
$r = $d.Content
$r.Find.Execute("1mmol above)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackPos = $r.Start
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $d.Range($goBackPos, $goBackPos)) | Out-Null

# ------------------------------------------------------------------
# 3) "onset of blood lactate accumulation (OBLA)" -> "aerobic threshold"
#    ("aerobic" becomes bold+underline, " threshold" stays bold only)
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("onset of blood lactate accumulation (OBLA)", $true, $false, $false, $false, $false, $true, 1, $false, "aerobic threshold", 2) | Out-Null
$sub = $d.Range($r.Start, $r.Start + 7)
$sub.Font.Underline = 1

# ------------------------------------------------------------------
# 4) "...we will feel that we are working, definitely above a tempo
#    pace." -> "...we will feel that we are working, but no more than
#    a tempo pace."
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("we will feel that we are working, definitely above a tempo pace", $true, $false, $false, $false, $false, $true, 1, $false, "we will feel that we are working, but no more than a tempo pace", 2) | Out-Null

# ------------------------------------------------------------------
# 5) "...called "LT2" or the maximal lactate steady state (MLSS) and
#    is closely related to..." -> "...called "LT2" or the anaerobic
#    threshold and more technically the maximal lactate steady state
#    (MLSS). The intensity / power at this point is closely related
#    to..."
#    ("anaerobic" becomes bold+underline, " threshold " stays bold
#    only, the "." right after MLSS stays bold)
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("maximal lactate steady state (MLSS)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$ins = $d.Range($r.Start, $r.Start)
$ins.InsertBefore("anaerobic threshold and more technically the ")
$anaerobic = $d.Range($ins.Start, $ins.Start + 9)
$anaerobic.Font.Bold = 1
$anaerobic.Font.Underline = 1
$thresh = $d.Range($ins.Start + 9, $ins.Start + 20)
$thresh.Font.Bold = 1

$r = $d.Content
$r.Find.Execute(" and is closely related to", $true, $false, $false, $false, $false, $true, 1, $false, ". The intensity / power at this point is closely related to", 2) | Out-Null
$dot = $d.Range($r.Start, $r.Start + 1)
$dot.Font.Bold = 1

# ------------------------------------------------------------------
# 6) " these are the destinations for that " - merge runs (bookmark
#    that used to sit here has been moved earlier in the document)
# ------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute(" these are the destinations for that ", $true, $false, $false, $false, $false, $true, 1, $false, " these are the destinations for that ", 2) | Out-Null
